# Refresh the crypto price/volume snapshot table with the latest scrape
# (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column ("D") holds values that look numeric (e.g. "1.010",
# "0.00000000370", "30.546.99") but must stay literal text, exactly as
# scraped, so force each touched cell to Text format before writing it.
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D8", "D9", "D10", "D11", "D13", "D15", "D16", "D18", "D19", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.546.99"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "2.135.06"
$ws.Range("E3").Value = "  +1.70%  "
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  +0.51%  "
$ws.Range("D5").Value = "352.03"
$ws.Range("E5").Value = "  +5.36%  "
$ws.Range("D6").Value = "1.009"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("E7").Value = "  +1.10%  "
$ws.Range("D8").Value = "0.4547"
$ws.Range("E8").Value = "  +0.63%  "
$ws.Range("D9").Value = "53.53"
$ws.Range("E9").Value = "  -0.65%  "
$ws.Range("D10").Value = "0.09168"
$ws.Range("E10").Value = "  +3.45%  "
$ws.Range("D11").Value = "1.183"
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("E12").Value = "  +5.01%  "
$ws.Range("D13").Value = "2.138.58"
$ws.Range("E13").Value = "  +2.03%  "
$ws.Range("E14").Value = "  +1.51%  "
$ws.Range("D15").Value = "8.194"
$ws.Range("E15").Value = "  +2.39%  "
$ws.Range("D16").Value = "102.17"
$ws.Range("E16").Value = "  +5.65%  "
$ws.Range("E17").Value = "  +2.86%  "
$ws.Range("D18").Value = "1.010"
$ws.Range("E18").Value = "  +0.49%  "
$ws.Range("D19").Value = "0.06723"
$ws.Range("E19").Value = "  +1.45%  "
$ws.Range("D20").Value = "20.30"
$ws.Range("E20").Value = "  +6.05%  "
$ws.Range("E21").Value = "  +0.44%  "
$ws.Range("D22").Value = "6.372"
$ws.Range("E22").Value = "  +1.78%  "
$ws.Range("D23").Value = "30.661.86"
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("D24").Value = "12.90"
$ws.Range("E24").Value = "  +4.91%  "
$ws.Range("D25").Value = "2.395"
$ws.Range("E25").Value = "  +2.57%  "
$ws.Range("D26").Value = "2.392.57"
$ws.Range("E26").Value = "  +2.20%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "2.661"
$ws.Range("E27").Value = "  +5.99%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "22.47"
$ws.Range("E28").Value = "  +1.58%  "
$ws.Range("D29").Value = "164.80"
$ws.Range("E29").Value = "  +1.41%  "
$ws.Range("D30").Value = "136.15"
$ws.Range("E30").Value = "  +2.53%  "
$ws.Range("E31").Value = "  +1.99%  "
$ws.Range("D32").Value = "0.1082"
$ws.Range("E32").Value = "  +1.61%  "
$ws.Range("D33").Value = "1.678"
$ws.Range("E33").Value = "  +1.95%  "
$ws.Range("D34").Value = "6.392"
$ws.Range("E34").Value = "  +0.59%  "
$ws.Range("D35").Value = "4.040"
$ws.Range("E35").Value = "  +2.33%  "
$ws.Range("D36").Value = "6.153"
$ws.Range("E36").Value = "  +5.86%  "
$ws.Range("E37").Value = "  +1.00%  "
$ws.Range("D38").Value = "0.02648"
$ws.Range("E38").Value = "  +3.01%  "
$ws.Range("D39").Value = "0.06971"
$ws.Range("E39").Value = "  +2.07%  "
$ws.Range("D40").Value = "0.2336"
$ws.Range("E40").Value = "  +2.03%  "
$ws.Range("D41").Value = "12.75"
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("D42").Value = "0.6987"
$ws.Range("E42").Value = "  +1.98%  "
$ws.Range("D43").Value = "1.281"
$ws.Range("E43").Value = "  +3.22%  "
$ws.Range("D44").Value = "14.77"
$ws.Range("E44").Value = "  +5.64%  "
$ws.Range("D45").Value = "2.360"
$ws.Range("E45").Value = "  +1.71%  "
$ws.Range("D46").Value = "0.6498"
$ws.Range("E46").Value = "  +2.61%  "
$ws.Range("D47").Value = "0.00000000370"
$ws.Range("E47").Value = "  +6.75%  "
$ws.Range("D48").Value = "3.754"
$ws.Range("E48").Value = "  +2.55%  "
$ws.Range("D49").Value = "1.252"
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("D50").Value = "84.12"
$ws.Range("E50").Value = "  +1.65%  "
$ws.Range("D51").Value = "0.07307"
$ws.Range("E51").Value = "  +2.69%  "
